$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remember the old content of rows 12-15 (these shift down to 16-19) ---
$c12 = $ws.Range("C12").Value2
$c13 = $ws.Range("C13").Value2
$c14 = $ws.Range("C14").Value2
$c15 = $ws.Range("C15").Value2

# --- Clear the old row 12-15 cells (they will be rewritten 4 rows lower) ---
$ws.Range("C12").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("C15").ClearContents()

# --- Rewrite that content 4 rows further down (row 12->16, 13->17, 14->18, 15->19) ---
$ws.Range("C16").Value2 = $c12
$ws.Range("C17").Value2 = $c13
$ws.Range("C18").Value2 = $c14
$ws.Range("C19").Value2 = $c15

# --- Existing "Interface" time entry on row 11: time changed from 15 to 75 ---
$ws.Range("C11").Value2 = 75

# --- New "Databas" activity block inserted at row 15 ---
$ws.Range("A15").Value2 = "Databas"

# Give the new date cell (B15) the same date number-format style as B11/B10 etc.
$ws.Range("B11").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("B15").Value2 = 44951

$ws.Range("C15").Value2 = 15

# --- New marker cell next to the Interface entry: a single space, stored as text ---
$ws.Range("F11").Value2 = " "

# --- Selection, as last left by the author ---
[void]$ws.Range("M6").Select()
